$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# M4:N5 array formula {=1,2;1,2}
$ws.Range("M4:N5").FormulaArray = "={1,2;1,2}"

# K6 regular formula (not array) containing the array literal text
$ws.Range("K6").Formula = "={1,2;1,2}"

# L6:M7 array formula K6+K5
$ws.Range("L6:M7").FormulaArray = "=K6+K5"

$ws.Range("L9").Select()
